$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.319.75'
$ws.Range("E2").Value = '  -0.15%  '

$ws.Range("D3").Value = '1.874.63'
$ws.Range("E3").Value = '  -0.07%  '

$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7098'
$ws.Range("E5").Value = '  -0.56%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.18'
$ws.Range("E6").Value = '  +0.05%  '

$ws.Range("E7").Value = '  +0.11%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07785'
$ws.Range("E8").Value = '  +0.55%  '

$ws.Range("E9").Value = '  -0.28%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '25.10'
$ws.Range("E10").Value = '  +0.90%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08419'
$ws.Range("E11").Value = '  +0.43%  '

$ws.Range("D12").Value = '1.868.60'
$ws.Range("E12").Value = '  -0.84%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.243'
$ws.Range("E13").Value = '  +0.04%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7172'
$ws.Range("E14").Value = '  +0.17%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.20'
$ws.Range("E15").Value = '  -0.17%  '

$ws.Range("D16").Value = '29.323.73'
$ws.Range("E16").Value = '  -0.06%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.090'
$ws.Range("E17").Value = '  +2.03%  '

$ws.Range("E18").Value = '  -0.13%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '240.70'
$ws.Range("E19").Value = '  -1.59%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.23'
$ws.Range("E20").Value = '  +0.27%  '

$ws.Range("D21").Value = '2.116.33'
$ws.Range("E21").Value = '  -0.87%  '

$ws.Range("E22").Value = '  +0.11%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.758'
$ws.Range("E23").Value = '  -1.96%  '

$ws.Range("E24").Value = '  +0.17%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1591'
$ws.Range("E25").Value = '  -1.85%  '

$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.45'
$ws.Range("E26").Value = '  -0.74%  '

$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.035'
$ws.Range("E27").Value = '  +0.12%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.51'
$ws.Range("E28").Value = '  -0.22%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.507'
$ws.Range("E29").Value = '  -0.07%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.410'
$ws.Range("E30").Value = '  -0.17%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.325'
$ws.Range("E31").Value = '  +0.43%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.287'
$ws.Range("E32").Value = '  -0.99%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05378'
$ws.Range("E33").Value = '  +3.15%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.946'
$ws.Range("E34").Value = '  +1.06%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7532'
$ws.Range("E35").Value = '  -2.61%  '

$ws.Range("E36").Value = '  +0.21%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.686'
$ws.Range("E37").Value = '  +0.11%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01888'
$ws.Range("E38").Value = '  +1.11%  '

$ws.Range("D39").Value = '1.239.13'
$ws.Range("E39").Value = '  +6.78%  '

$ws.Range("E40").Value = '  +0.68%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.482'
$ws.Range("E41").Value = '  +1.51%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8934'
$ws.Range("E42").Value = '  +0.14%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '72.40'
$ws.Range("E43").Value = '  -1.58%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '108.62'
$ws.Range("E44").Value = '  +4.57%  '

$ws.Range("E45").Value = '  +0.14%  '

$ws.Range("D46").Value = '2.018.05'
$ws.Range("E46").Value = '  -0.20%  '

$ws.Range("E47").Value = '  +6.75%  '

$ws.Range("E48").Value = '  -0.48%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.5199'
$ws.Range("E49").Value = '  +0.15%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.450'
$ws.Range("E50").Value = '  +0.38%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4339'
$ws.Range("E51").Value = '  +0.55%  '
